$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Week 46" column (AU) header
$ws.Range("AU1").Value = "Week 46"

# Day-after inputs for Week 46 (column AU)
$ws.Range("AU6").Value = 5.5
$ws.Range("AU7").Value = 5.5
$ws.Range("AU8").Value = 4.25
$ws.Range("AU9").Value = 1.5
$ws.Range("AU10").Value = 6

# Restore the active selection (shifted one column left of the new last
# column, matching the author's saved cursor position for this edit).
$ws.Range("AP10").Select()
